$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 0.15
$ws.Range("N3").Value = 0.1398479931319728
$ws.Range("N4").Value = 0.007922667858385166
$ws.Range("N6").Value = 0.001
$ws.Range("N7").Value = 0.001000000000000005
$ws.Range("N8").Value = 0.15
$ws.Range("N9").Value = 0.001
$ws.Range("N11").Value = 0.07296069780610072
$ws.Range("N12").Value = 0.06806522584086927
$ws.Range("N13").Value = 0.03061350869190371
$ws.Range("N14").Value = 0.15
$ws.Range("N15").Value = 0.001
$ws.Range("N16").Value = 0.07526094837697749
$ws.Range("N17").Value = 0.01022751055084323
$ws.Range("N18").Value = 0.001
$ws.Range("N20").Value = 0.1361014477429478
$ws.Range("N21").Value = 0.001000000000000001
$ws.Range("B22").Value = -0.01873928690977155
$ws.Range("C22").Value = 0.08335350925606616
$ws.Range("D22").Value = 0.09346953300218509
$ws.Range("E22").Value = -0.1328983601556065
$ws.Range("F22").Value = -0.03016696668154265
$ws.Range("G22").Value = -0.01654894083654105
$ws.Range("H22").Value = 0.02576982480353973
$ws.Range("I22").Value = 0.07630795323203993
$ws.Range("J22").Value = 0.03505692427426015
$ws.Range("K22").Value = 0.038204404226788
$ws.Range("L22").Value = -0.02202903048723463
$ws.Range("M22").Value = -0.03129112320296311
$ws.Range("B23").Value = 0.9814352018953287
$ws.Range("C23").Value = 1.086925979034573
$ws.Range("D23").Value = 1.097977150764249
$ws.Range("E23").Value = 0.8755540787667234
$ws.Range("F23").Value = 0.9702835150043794
$ws.Range("G23").Value = 0.9835872406303389
$ws.Range("H23").Value = 1.026104737429997
$ws.Range("I23").Value = 1.07929489533631
$ws.Range("J23").Value = 1.035678662377885
$ws.Range("K23").Value = 1.038943575634134
$ws.Range("L23").Value = 0.9782118366727653
$ws.Range("M23").Value = 0.9691933773208566
$ws.Range("N23").Value = 1.130109998234451

Write-Output "applied edits"
